# Rebuild the "Model_Performance_Summary" table with the new
# Model_Family / Engine / Split_Type columns (inserted after "Model" and
# before "Source"), and re-derive Source into Split_Type values
# (TS_CV / N/A) while Model_Family holds GARCH / NF-GARCH and Engine is
# always N/A for this consolidated export. Row order also changes
# slightly because results are now grouped by split type.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$arr = New-Object 'object[,]' 10,10
$arr[0,0] = "Model"
$arr[0,1] = "Model_Family"
$arr[0,2] = "Engine"
$arr[0,3] = "Split_Type"
$arr[0,4] = "Source"
$arr[0,5] = "Avg_AIC"
$arr[0,6] = "Avg_BIC"
$arr[0,7] = "Avg_LogLik"
$arr[0,8] = "Avg_MSE"
$arr[0,9] = "Avg_MAE"
$arr[1,0] = "TGARCH"
$arr[1,1] = "GARCH"
$arr[1,2] = "N/A"
$arr[1,3] = "TS_CV"
$arr[1,4] = "Time_Series_CV"
$arr[1,5] = -6.23674331462475
$arr[1,6] = -6.17773880124684
$arr[1,7] = 1566.18582865619
$arr[1,8] = 0.000265231788437081
$arr[1,9] = 0.00997959312049119
$arr[2,0] = "eGARCH"
$arr[2,1] = "GARCH"
$arr[2,2] = "N/A"
$arr[2,3] = "TS_CV"
$arr[2,4] = "Time_Series_CV"
$arr[2,5] = -6.16499588686336
$arr[2,6] = -6.10599137348545
$arr[2,7] = 1548.24897171584
$arr[2,8] = 0.000265255247246037
$arr[2,9] = 0.00998057275133209
$arr[3,0] = "eGARCH"
$arr[3,1] = "NF-GARCH"
$arr[3,2] = "N/A"
$arr[3,3] = "N/A"
$arr[3,4] = "NF-GARCH"
$arr[3,5] = 27401.3105107422
$arr[3,6] = 27439.8014735373
$arr[3,7] = -13694.6552553711
$arr[3,8] = 0
$arr[3,9] = 0
$arr[4,0] = "fGARCH"
$arr[4,1] = "NF-GARCH"
$arr[4,2] = "N/A"
$arr[4,3] = "N/A"
$arr[4,4] = "NF-GARCH"
$arr[4,5] = -28273.7513163331
$arr[4,6] = -28235.260353538
$arr[4,7] = 14142.8756581666
$arr[4,8] = 0
$arr[4,9] = 0
$arr[5,0] = "gjrGARCH"
$arr[5,1] = "GARCH"
$arr[5,2] = "N/A"
$arr[5,3] = "TS_CV"
$arr[5,4] = "Time_Series_CV"
$arr[5,5] = -6.23529221703763
$arr[5,6] = -6.17628770365972
$arr[5,7] = 1565.82305425941
$arr[5,8] = 0.000264427773901888
$arr[5,9] = 0.0099650836905736
$arr[6,0] = "gjrGARCH"
$arr[6,1] = "NF-GARCH"
$arr[6,2] = "N/A"
$arr[6,3] = "N/A"
$arr[6,4] = "NF-GARCH"
$arr[6,5] = -28244.7567895128
$arr[6,6] = -28206.2658267177
$arr[6,7] = 14128.3783947564
$arr[6,8] = 0
$arr[6,9] = 0
$arr[7,0] = "sGARCH"
$arr[7,1] = "NF-GARCH"
$arr[7,2] = "N/A"
$arr[7,3] = "N/A"
$arr[7,4] = "NF-GARCH"
$arr[7,5] = -27993.3295905363
$arr[7,6] = -27964.46136844
$arr[7,7] = 14001.1647952681
$arr[7,8] = 0
$arr[7,9] = 0
$arr[8,0] = "sGARCH_norm"
$arr[8,1] = "GARCH"
$arr[8,2] = "N/A"
$arr[8,3] = "TS_CV"
$arr[8,4] = "Time_Series_CV"
$arr[8,5] = -6.14008516085673
$arr[8,6] = -6.10636829606935
$arr[8,7] = 1539.02129021418
$arr[8,8] = 0.00026524306076551
$arr[8,9] = 0.00997754931432829
$arr[9,0] = "sGARCH_sstd"
$arr[9,1] = "GARCH"
$arr[9,2] = "N/A"
$arr[9,3] = "TS_CV"
$arr[9,4] = "Time_Series_CV"
$arr[9,5] = -6.22901892609995
$arr[9,6] = -6.17844362891889
$arr[9,7] = 1563.25473152499
$arr[9,8] = 0.000265218411076707
$arr[9,9] = 0.00997773042969076
$ws.Range("A1:J10").Value = $arr
